$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to be stored as text so that
# numeric-looking values (e.g. "12.29") are not auto-converted to numbers,
# matching the original inline-string cell content.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.389.90'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.841.57'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '239.17'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '0.6247'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '0.07387'
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '0.2889'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = '24.86'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").Value = '0.07714'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.841.31'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '4.961'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '0.00001025'
$ws.Range("E15").Value = '  -1.84%  '
$ws.Range("D16").Value = '81.75'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '6.279'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").Value = '29.346.03'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '234.40'
$ws.Range("E19").Value = '  +2.81%  '
$ws.Range("D20").Value = '12.29'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '7.295'
$ws.Range("E22").Value = '  -2.76%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '157.18'
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").Value = '8.473'
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("E26").Value = '  -1.93%  '
$ws.Range("D27").Value = '17.30'
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("D28").Value = '0.07293'
$ws.Range("E28").Value = '  +13.50%  '
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  +5.39%  '
$ws.Range("D30").Value = '1.475'
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").Value = '4.041'
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").Value = '4.028'
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '1.157'
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '1.817'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("D35").Value = '0.7035'
$ws.Range("E35").Value = '  +1.26%  '
$ws.Range("D36").Value = '2.574'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").Value = '0.01830'
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").Value = '2.785'
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("D39").Value = '1.234.36'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").Value = '6.764'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").Value = '0.9501'
$ws.Range("E41").Value = '  +2.10%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '1.990.70'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").Value = '101.13'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = '65.29'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").Value = '0.00000000118'
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("D47").Value = '6.972'
$ws.Range("E47").Value = '  -1.40%  '
$ws.Range("D48").Value = '1.698'
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("D49").Value = '8.876'
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1132'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.3882'
$ws.Range("E51").Value = '  -1.45%  '
